$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update the acquisition-timestamp column (A) for all data rows (2-11)
# to reflect the latest scrape run appended at 2025-11-27 12:50:22.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-11-27 12:50:22"
}
